# Update "想去人数" (want-to-go count) values in column F across sheets.
$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F5").Value = 1034
$ws1.Range("F6").Value = 207
$ws1.Range("F8").Value = 112
$ws1.Range("F9").Value = 563
$ws1.Range("F10").Value = 578
$ws1.Range("F13").Value = 146

# Sheet 2: 演出 (Show)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 8

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6270
$ws3.Range("F4").Value = 1890

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6270
$ws4.Range("F4").Value = 1890
$ws4.Range("F9").Value = 8
$ws4.Range("F14").Value = 1034
$ws4.Range("F16").Value = 207
$ws4.Range("F20").Value = 112
$ws4.Range("F21").Value = 563
$ws4.Range("F23").Value = 578
$ws4.Range("F28").Value = 146
